$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Win probabilities" sheet: the old sheet had a helper header row (row 1:
# team numbers 1..16) and a helper header column (col A: team numbers 1..16)
# around a 16x16 symmetric win-probability matrix. Both the header row and
# header column are removed, leaving just the matrix in A1:P16.
# ---------------------------------------------------------------------------
$wsWin = $wb.Worksheets.Item("Win probabilities")
$wsWin.Rows.Item(1).Delete()
$wsWin.Columns.Item(1).Delete()

# ---------------------------------------------------------------------------
# "Tournament" sheet: add the bracket's winner-calculation helper formulas
# that reference the (now shifted) "Win probabilities" matrix, and rework the
# two existing simulation formulas (C2 / C6) to use the new local helper
# cells instead of reaching directly into 'Win probabilities'.
# ---------------------------------------------------------------------------
$wsT = $wb.Worksheets.Item("Tournament")

# Round 1 (top bracket, teams 1-4) ------------------------------------------------
$wsT.Range("A2").Formula  = "='Win probabilities'!A2"
$wsT.Range("C2").Formula  = "=IF(RAND()<A2,A1,A3)"

$wsT.Range("A6").Formula  = "='Win probabilities'!C4"
$wsT.Range("C6").Formula  = "=IF(RAND()<A6,A5,A7)"

$wsT.Range("C1").Formula  = "=AND(C2=1,C6=3)"
$wsT.Range("D1").Formula  = "=IF(C1=TRUE,'Win probabilities'!A3)"
$wsT.Range("C3").Formula  = "=AND(C2=2,C6=3)"
$wsT.Range("D3").Formula  = "=IF(C3=TRUE,'Win probabilities'!B3)"
$wsT.Range("C5").Formula  = "=AND(C2=1,C6=4)"
$wsT.Range("D5").Formula  = "=IF(C5=TRUE,'Win probabilities'!A4)"
$wsT.Range("C7").Formula  = "=AND(C2=2,C6=4)"
$wsT.Range("D7").Formula  = "=IF(C7=TRUE,'Win probabilities'!B4)"

$wsT.Range("E4").Formula  = "=IF(RAND()<SUM(D1:D7),C2,C6)"

# Round 1 (bottom bracket, teams 5-8) ---------------------------------------------
$wsT.Range("A10").Formula = "='Win probabilities'!E6"
$wsT.Range("C11").Formula = "=IF(RAND()<A10,A9,A12)"

$wsT.Range("A15").Formula = "='Win probabilities'!G8"
$wsT.Range("C15").Formula = "=IF(RAND()<A15,A14,A16)"

$wsT.Range("C10").Formula = "=AND(C11=5,C15=7)"
$wsT.Range("D10").Formula = "=IF(C10=TRUE,'Win probabilities'!E7)"
$wsT.Range("C12").Formula = "=AND(C11=6,C15=7)"
$wsT.Range("D12").Formula = "=IF(C12=TRUE,'Win probabilities'!F7)"
$wsT.Range("C14").Formula = "=AND(C11=5,C15=8)"
$wsT.Range("D14").Formula = "=IF(C14=TRUE,'Win probabilities'!E8)"
$wsT.Range("C16").Formula = "=AND(C11=6,C15=8)"
$wsT.Range("D16").Formula = "=IF(C16=TRUE,'Win probabilities'!F8)"

$wsT.Range("E13").Formula = "=IF(RAND()<SUM(D10:D16),C11,C15)"

# New team-number labels on the right-hand bracket column (K1..K16), and the
# newly-filled-in team numbers on the left/middle column (A9, A12, A14, A16).
$wsT.Range("K1").Value  = 9
$wsT.Range("K3").Value  = 10
$wsT.Range("K5").Value  = 11
$wsT.Range("K7").Value  = 12
$wsT.Range("K9").Value  = 13
$wsT.Range("K12").Value = 14
$wsT.Range("K14").Value = 15
$wsT.Range("K16").Value = 16
$wsT.Range("A9").Value  = 5
$wsT.Range("A12").Value = 6
$wsT.Range("A14").Value = 7
$wsT.Range("A16").Value = 8

# I2 previously used a "no top border" box style; it now matches the plain
# bordered-box style used by the other helper cells in column I (I6/I11/I15).
$wsT.Range("I6").Copy()
$wsT.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the original selections on each sheet.
$wsWin.Range("D17").Select()
$wsT.Activate()
$wsT.Range("D1").Select()
